# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" worksheets to reflect refreshed scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 75
    4  = 12183
    5  = 4485
    7  = 61
    10 = 2601
    11 = 1128
    12 = 196
    13 = 68
    14 = 5285
    15 = 66
    16 = 206
    17 = 551
    18 = 11454
    19 = 11521
    21 = 59
    25 = 27
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# Sheet "全部类型" (All types) - row => new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    3  = 75
    4  = 12183
    5  = 4485
    7  = 61
    10 = 2601
    12 = 1128
    13 = 196
    14 = 68
    15 = 5285
    16 = 66
    17 = 206
    18 = 551
    19 = 11454
    20 = 11521
    22 = 59
    26 = 27
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}

$wb.Save()
